# Update the "Sales" values in column B (rows 2-61) of Sheet1
# per the diff: B2 through B61 get new values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValues = @(
    187523.4001659137,
    105562.7767840706,
    205210.7663225566,
    228811.9890069297,
    120416.8967446932,
    139393.1624332218,
    275159.8463265749,
    168214.0150567523,
    151705.3271741188,
    187837.3638573977,
    288816.8084868879,
    158830.0785836617,
    224828.4016010703,
    255284.2746951019,
    170889.9193569951,
    257274.5142980206,
    200845.6489483423,
    222862.5790626067,
    202753.0500662491,
    191585.3469876064,
    101581.0049861069,
    245219.3246358919,
    204476.4938604359,
    139128.489028873,
    168589.5211168695,
    208772.2724329935,
    243979.3617480665,
    123049.5793452091,
    219595.7670896158,
    169857.7399855685,
    208513.6981587364,
    164190.1454452684,
    211072.4208534042,
    140464.336236545,
    273378.6368979879,
    127984.1064938287,
    209152.5213140719,
    173929.8702813017,
    285216.3332170502,
    250071.5055673061,
    127919.635303659,
    235432.3982864047,
    231454.0649600415,
    250359.6177525136,
    224770.1322014902,
    174487.7126038608,
    114878.5845531393,
    263378.6544322563,
    267403.6696283267,
    209775.0199613006,
    150017.7756329985,
    183036.4210796841,
    254095.5262188889,
    262545.7776712077,
    264435.4218551953,
    226895.2446206884,
    229183.2462271144,
    163382.8058268771,
    122952.5647260233,
    260619.9081770256
)

$startRow = 2
for ($i = 0; $i -lt $newValues.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 2).Value = $newValues[$i]
}
